$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 15.13657353076129
$ws.Range("C2").Value = 11.17067940586771
$ws.Range("D2").Value = 3.851690780264482
$ws.Range("F2").Value = 18.77712301612766
$ws.Range("G2").Value = 19.39772137241023
$ws.Range("H2").Value = 11.83593624816265
$ws.Range("I2").Value = 15.84101858613838
# Row 3
$ws.Range("B3").Value = 14.30321865981795
$ws.Range("C3").Value = 10.50151679589502
$ws.Range("D3").Value = 3.78741438666288
$ws.Range("F3").Value = 18.87249940632003
$ws.Range("G3").Value = 19.5284803802863
$ws.Range("H3").Value = 11.92511203841327
$ws.Range("I3").Value = 16.04036731711333
# Row 4
$ws.Range("B4").Value = 13.76561161901914
$ws.Range("C4").Value = 10.06644530425696
$ws.Range("D4").Value = 3.747447471700184
$ws.Range("F4").Value = 18.94204487394068
$ws.Range("G4").Value = 19.6269469956593
$ws.Range("H4").Value = 11.9837130206944
$ws.Range("I4").Value = 16.16902746109914
# Row 5
$ws.Range("B5").Value = 13.54018228269856
$ws.Range("C5").Value = 9.883105675976187
$ws.Range("D5").Value = 3.73105070080506
$ws.Range("F5").Value = 18.97311741772196
$ws.Range("G5").Value = 19.67156606130252
$ws.Range("H5").Value = 12.00855751083072
$ws.Range("I5").Value = 16.2230339394446
# Row 6
$ws.Range("B6").Value = 13.50237200546249
$ws.Range("C6").Value = 9.852298667727055
$ws.Range("D6").Value = 3.728321866988986
$ws.Range("F6").Value = 18.97844110796926
$ws.Range("G6").Value = 19.67924402254084
$ws.Range("H6").Value = 12.01274105895263
$ws.Range("I6").Value = 16.23209696159923
# Row 7
$ws.Range("B7").Value = 13.76259686301858
$ws.Range("C7").Value = 10.06399713154697
$ws.Range("D7").Value = 3.747226762723684
$ws.Range("F7").Value = 18.94245290709963
$ws.Range("G7").Value = 19.62753065833657
$ws.Range("H7").Value = 11.98404418265954
$ws.Range("I7").Value = 16.16974942513139
# Row 8
$ws.Range("B8").Value = 14.85470406809883
$ws.Range("C8").Value = 10.94499662481838
$ws.Range("D8").Value = 3.829641737158833
$ws.Range("F8").Value = 18.80771161580034
$ws.Range("G8").Value = 19.43898805756712
$ws.Range("H8").Value = 11.86588375199377
$ws.Range("I8").Value = 15.9084564115531
# Row 9
$ws.Range("B9").Value = 16.78481528197226
$ws.Range("C9").Value = 12.47934854740182
$ws.Range("D9").Value = 3.986602071055713
$ws.Range("F9").Value = 18.63188793262519
$ws.Range("G9").Value = 19.21675199175748
$ws.Range("H9").Value = 11.66483823315994
$ws.Range("I9").Value = 15.44561306761325
# Row 10
$ws.Range("B10").Value = 18.06795970058275
$ws.Range("C10").Value = 13.48777835697966
$ws.Range("D10").Value = 4.098210854966283
$ws.Range("F10").Value = 18.55817992912533
$ws.Range("G10").Value = 19.14738255427066
$ws.Range("H10").Value = 11.53602937671541
$ws.Range("I10").Value = 15.13561511307443
# Row 11
$ws.Range("B11").Value = 18.62155576339551
$ws.Range("C11").Value = 13.9206507808591
$ws.Range("D11").Value = 4.148006669693518
$ws.Range("F11").Value = 18.5369856754886
$ws.Range("G11").Value = 19.13692110414199
$ws.Range("H11").Value = 11.48158102848609
$ws.Range("I11").Value = 15.00108477095148
# Row 12
$ws.Range("B12").Value = 18.82680754233744
$ws.Range("C12").Value = 14.08084737017495
$ws.Range("D12").Value = 4.166709272799856
$ws.Range("F12").Value = 18.53075430042479
$ws.Range("G12").Value = 19.13603928129569
$ws.Range("H12").Value = 11.46156338235683
$ws.Range("I12").Value = 14.95107301294702
# Row 13
$ws.Range("B13").Value = 18.78279854488817
$ws.Range("C13").Value = 14.04651175052567
$ws.Range("D13").Value = 4.162688391732809
$ws.Range("F13").Value = 18.53201623193355
$ws.Range("G13").Value = 19.13609158640968
$ws.Range("H13").Value = 11.46584775063095
$ws.Range("I13").Value = 14.96180250719972
# Row 14
$ws.Range("B14").Value = 18.63853005147424
$ws.Range("C14").Value = 13.93390486555145
$ws.Range("D14").Value = 4.149548506192589
$ws.Range("F14").Value = 18.53643694522144
$ws.Range("G14").Value = 19.13678660124024
$ws.Range("H14").Value = 11.47992209905442
$ws.Range("I14").Value = 14.99695160165079
# Row 15
$ws.Range("B15").Value = 18.54958928006608
$ws.Range("C15").Value = 13.86444505103866
$ws.Range("D15").Value = 4.141479498556861
$ws.Range("F15").Value = 18.53937901201568
$ws.Range("G15").Value = 19.13761460578099
$ws.Range("H15").Value = 11.48862141258787
$ws.Range("I15").Value = 15.01860277433099
# Row 16
$ws.Range("B16").Value = 18.03117013658767
$ws.Range("C16").Value = 13.4589682536296
$ws.Range("D16").Value = 4.094935761982139
$ws.Range("F16").Value = 18.55981506884858
$ws.Range("G16").Value = 19.14849487812871
$ws.Range("H16").Value = 11.53967153148595
$ws.Range("I16").Value = 15.1445374619281
# Row 17
$ws.Range("B17").Value = 17.70538239048791
$ws.Range("C17").Value = 13.20359264026464
$ws.Range("D17").Value = 4.066122835090214
$ws.Range("F17").Value = 18.57552688524897
$ws.Range("G17").Value = 19.16060823158278
$ws.Range("H17").Value = 11.57205444848863
$ws.Range("I17").Value = 15.22345546354519
# Row 18
$ws.Range("B18").Value = 17.51516578781625
$ws.Range("C18").Value = 13.05427248991933
$ws.Range("D18").Value = 4.049459349070616
$ws.Range("F18").Value = 18.58572394955862
$ws.Range("G18").Value = 19.16955885163669
$ws.Range("H18").Value = 11.59107035361101
$ws.Range("I18").Value = 15.26945787089329
# Row 19
$ws.Range("B19").Value = 17.45027712150766
$ws.Range("C19").Value = 13.0032972240334
$ws.Range("D19").Value = 4.043802170225386
$ws.Range("F19").Value = 18.58937507028557
$ws.Range("G19").Value = 19.17292846401643
$ws.Range("H19").Value = 11.59757568302677
$ws.Range("I19").Value = 15.28513845108737
# Row 20
$ws.Range("B20").Value = 17.74035658809655
$ws.Range("C20").Value = 13.23102984300602
$ws.Range("D20").Value = 4.069199545530646
$ws.Range("F20").Value = 18.57373411828462
$ws.Range("G20").Value = 19.15911311280274
$ws.Range("H20").Value = 11.56856682579081
$ws.Range("I20").Value = 15.2149912950307
# Row 21
$ws.Range("B21").Value = 18.68102451841392
$ws.Range("C21").Value = 13.96708127910273
$ws.Range("D21").Value = 4.153412293922214
$ws.Range("F21").Value = 18.53508962368056
$ws.Range("G21").Value = 19.13649855224072
$ws.Range("H21").Value = 11.47577178180101
$ws.Range("I21").Value = 14.98660217579198
# Row 22
$ws.Range("B22").Value = 19.27024160115353
$ws.Range("C22").Value = 14.42643517709242
$ws.Range("D22").Value = 4.207546782005295
$ws.Range("F22").Value = 18.52030014465709
$ws.Range("G22").Value = 19.13968447206701
$ws.Range("H22").Value = 11.41862867664924
$ws.Range("I22").Value = 14.84276847890469
# Row 23
$ws.Range("B23").Value = 18.9581187564113
$ws.Range("C23").Value = 14.18325483030501
$ws.Range("D23").Value = 4.178741184967524
$ws.Range("F23").Value = 18.52722974733516
$ws.Range("G23").Value = 19.13632727894406
$ws.Range("H23").Value = 11.44880495040506
$ws.Range("I23").Value = 14.91903855714716
# Row 24
$ws.Range("B24").Value = 17.72455383471403
$ws.Range("C24").Value = 13.21863327452379
$ws.Range("D24").Value = 4.067808870641834
$ws.Range("F24").Value = 18.57454100358317
$ws.Range("G24").Value = 19.15978287182916
$ws.Range("H24").Value = 11.57014233887668
$ws.Range("I24").Value = 15.21881597865865
# Row 25
$ws.Range("B25").Value = 16.2860150313629
$ws.Range("C25").Value = 12.08509275985459
$ws.Range("D25").Value = 3.944730889141612
$ws.Range("F25").Value = 18.66981803648606
$ws.Range("G25").Value = 19.26064501161281
$ws.Range("H25").Value = 11.71592269490779
$ws.Range("I25").Value = 15.56553626496413
